$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-07-30 Sunday" "2023-07-31 Monday"

Replace-Text "61÷2=30, 1" "18÷9=2, 0"
Replace-Text "41÷4=10, 1" "91÷5=18, 1"
Replace-Text "86÷2=43, 0" "14÷5=2, 4"
Replace-Text "44÷5=8, 4" "84÷8=10, 4"
Replace-Text "15÷5=3, 0" "53÷9=5, 8"

Replace-Text "98÷3=32, 2" "65÷9=7, 2"
Replace-Text "94÷7=13, 3" "23÷7=3, 2"
Replace-Text "27÷6=4, 3" "30÷6=5, 0"
Replace-Text "77÷2=38, 1" "83÷6=13, 5"
Replace-Text "70÷3=23, 1" "44÷4=11, 0"

Replace-Text "91÷4=22, 3" "13÷3=4, 1"
Replace-Text "59÷4=14, 3" "42÷6=7, 0"
Replace-Text "65÷2=32, 1" "25÷8=3, 1"
Replace-Text "74÷9=8, 2" "76÷5=15, 1"
Replace-Text "29÷5=5, 4" "94÷9=10, 4"

Replace-Text "89÷3=29, 2" "87÷4=21, 3"
Replace-Text "12÷6=2, 0" "31÷8=3, 7"
Replace-Text "10÷7=1, 3" "97÷7=13, 6"
Replace-Text "47÷9=5, 2" "84÷2=42, 0"
Replace-Text "59÷5=11, 4" "23÷5=4, 3"

Replace-Text "11÷3=3, 2" "96÷3=32, 0"
Replace-Text "60÷8=7, 4" "13÷9=1, 4"
Replace-Text "94÷8=11, 6" "64÷6=10, 4"
Replace-Text "58÷9=6, 4" "96÷7=13, 5"
Replace-Text "65÷6=10, 5" "66÷9=7, 3"
